# Corrected GS transmit power
$wb = $excel.ActiveWorkbook

# Update GS Transmit Power (GSTP, Input!C14) from 100 to 50
$inputSheet = $wb.Worksheets.Item("Input")
$inputSheet.Range("C14").Value = 50

# Make the "Input" sheet the active/selected sheet (tabSelected) and set its selection
$inputSheet.Activate()
$inputSheet.Range("J15").Select()

$wb.Save()
